# Auto-generated: updates currentAveragePrice / LevePrice / LeveProfit
# columns across sheets with refreshed market-board data (scheduled runner).
$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 0
$ws.Range("J18").Value = 0
$ws.Range("L18").Value = 0
$ws.Range("N18").ClearContents()
$ws.Range("H33").Value = 136.46153
$ws.Range("I33").Value = 135
$ws.Range("K33").Value = 135
$ws.Range("M33").Value = 94
$ws.Range("H34").Value = 2230.5557
$ws.Range("I34").Value = 2230.5557
$ws.Range("K34").Value = 2230.5557
$ws.Range("M34").Value = -2027.5557
$ws.Range("H36").Value = 2230.5557
$ws.Range("I36").Value = 2230.5557
$ws.Range("K36").Value = 2230.5557
$ws.Range("M36").Value = -1515.5557
$ws.Range("H106").Value = 34186.375
$ws.Range("I106").Value = 36249.332
$ws.Range("K106").Value = 36249.332
$ws.Range("M106").Value = -35618.332
$ws.Range("H129").Value = 2124.5454
$ws.Range("I129").Value = 895.1667
$ws.Range("J129").Value = 3599.8
$ws.Range("K129").Value = 2685.5001
$ws.Range("L129").Value = 10799.4
$ws.Range("M129").Value = 2314.4999
$ws.Range("N129").Value = -20799.4
$ws.Range("H132").Value = 1078.4783
$ws.Range("I132").Value = 1123.5238
$ws.Range("K132").Value = 3370.5714
$ws.Range("M132").Value = -840.5713999999998
$ws.Range("H137").Value = 1730.1765
$ws.Range("I137").Value = 1658
$ws.Range("K137").Value = 4974
$ws.Range("M137").Value = -2424
$ws.Range("H138").Value = 4598.773
$ws.Range("I138").Value = 3341.9355
$ws.Range("J138").Value = 5711.971
$ws.Range("K138").Value = 10025.8065
$ws.Range("L138").Value = 17135.913
$ws.Range("M138").Value = -4885.806500000001
$ws.Range("N138").Value = -27415.913

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 4144.9165
$ws.Range("I122").Value = 3923.9
$ws.Range("J122").Value = 5250
$ws.Range("K122").Value = 11771.7
$ws.Range("L122").Value = 15750
$ws.Range("M122").Value = -9321.700000000001
$ws.Range("N122").Value = -20650
$ws.Range("H132").Value = 1038.6666
$ws.Range("I132").Value = 1038.6666
$ws.Range("K132").Value = 3115.9998
$ws.Range("M132").Value = -585.9998000000001

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 4691.5
$ws.Range("I20").Value = 6300.8
$ws.Range("K20").Value = 6300.8
$ws.Range("M20").Value = -6053.8
$ws.Range("H99").Value = 2077.348
$ws.Range("I99").Value = 1531.5834
$ws.Range("K99").Value = 1531.5834
$ws.Range("M99").Value = -33.58339999999998

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 135.94118
$ws.Range("I7").Value = 78.84614999999999
$ws.Range("K7").Value = 78.84614999999999
$ws.Range("M7").Value = 34.15385000000001
$ws.Range("H12").Value = 1000
$ws.Range("I12").Value = 0
$ws.Range("J12").Value = 1000
$ws.Range("K12").Value = 0
$ws.Range("L12").Value = 1000
$ws.Range("M12").ClearContents()
$ws.Range("N12").Value = -1340
$ws.Range("H31").Value = 4111.75
$ws.Range("I31").Value = 3595.5881
$ws.Range("K31").Value = 3595.5881
$ws.Range("M31").Value = -3300.5881
$ws.Range("H34").Value = 4111.75
$ws.Range("I34").Value = 3595.5881
$ws.Range("K34").Value = 3595.5881
$ws.Range("M34").Value = -3393.5881
$ws.Range("H58").Value = 2375.7083
$ws.Range("I58").Value = 1292.1578
$ws.Range("K58").Value = 1292.1578
$ws.Range("M58").Value = -1089.1578
$ws.Range("H134").Value = 3637.5
$ws.Range("I134").Value = 3565.2
$ws.Range("K134").Value = 10695.6
$ws.Range("M134").Value = -8160.599999999999
$ws.Range("H136").Value = 2375.7083
$ws.Range("I136").Value = 1292.1578
$ws.Range("K136").Value = 3876.4734
$ws.Range("M136").Value = -1326.4734

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H82").Value = 18499.75
$ws.Range("J82").Value = 18499.75
$ws.Range("L82").Value = 55499.25
$ws.Range("N82").Value = -56311.25
$ws.Range("H85").Value = 18499.75
$ws.Range("J85").Value = 18499.75
$ws.Range("L85").Value = 55499.25
$ws.Range("N85").Value = -58307.25

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H59").Value = 0
$ws.Range("I59").Value = 0
$ws.Range("K59").Value = 0
$ws.Range("M59").ClearContents()
$ws.Range("H92").Value = 15499.75
$ws.Range("J92").Value = 12333
$ws.Range("L92").Value = 12333
$ws.Range("N92").Value = -16077
$ws.Range("H93").Value = 56744.75
$ws.Range("J93").Value = 56744.75
$ws.Range("L93").Value = 56744.75
$ws.Range("N93").Value = -60488.75
$ws.Range("H110").Value = 80000
$ws.Range("J110").Value = 80000
$ws.Range("L110").Value = 80000
$ws.Range("N110").Value = -88180
$ws.Range("H119").Value = 95000
$ws.Range("J119").Value = 95000
$ws.Range("L119").Value = 95000
$ws.Range("N119").Value = -104676
$ws.Range("H122").Value = 114403.445
$ws.Range("I122").Value = 3132.6
$ws.Range("J122").Value = 253492
$ws.Range("K122").Value = 9397.799999999999
$ws.Range("L122").Value = 760476
$ws.Range("M122").Value = -6947.799999999999
$ws.Range("N122").Value = -765376
$ws.Range("H123").Value = 30336.3
$ws.Range("J123").Value = 30336.3
$ws.Range("L123").Value = 30336.3
$ws.Range("N123").Value = -35236.3
$ws.Range("H132").Value = 2813.1333
$ws.Range("I132").Value = 2516.4167
$ws.Range("K132").Value = 7549.250100000001
$ws.Range("M132").Value = -5019.250100000001
$ws.Range("H133").Value = 143310.67
$ws.Range("J133").Value = 143310.67
$ws.Range("L133").Value = 143310.67
$ws.Range("N133").Value = -153430.67

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 5117
$ws.Range("I16").Value = 5419.4287
$ws.Range("K16").Value = 5419.4287
$ws.Range("M16").Value = -5249.4287
$ws.Range("H22").Value = 2402.9285
$ws.Range("I22").Value = 1738.6897
$ws.Range("J22").Value = 3884.6924
$ws.Range("K22").Value = 1738.6897
$ws.Range("L22").Value = 3884.6924
$ws.Range("M22").Value = -1443.6897
$ws.Range("N22").Value = -4474.6924
$ws.Range("H27").Value = 2402.9285
$ws.Range("I27").Value = 1738.6897
$ws.Range("J27").Value = 3884.6924
$ws.Range("K27").Value = 1738.6897
$ws.Range("L27").Value = 3884.6924
$ws.Range("M27").Value = -1631.6897
$ws.Range("N27").Value = -4098.6924
$ws.Range("H82").Value = 2264.1667
$ws.Range("I82").Value = 2264.1667
$ws.Range("K82").Value = 2264.1667
$ws.Range("M82").Value = -1903.1667
$ws.Range("H85").Value = 2264.1667
$ws.Range("I85").Value = 2264.1667
$ws.Range("K85").Value = 2264.1667
$ws.Range("M85").Value = -1016.1667
$ws.Range("H132").Value = 3799.8
$ws.Range("I132").Value = 3799.8
$ws.Range("K132").Value = 11399.4
$ws.Range("M132").Value = -8869.400000000001
$ws.Range("H133").Value = 60000
$ws.Range("J133").Value = 60000
$ws.Range("L133").Value = 60000
$ws.Range("N133").Value = -65060
$ws.Range("H136").Value = 3867.3333
$ws.Range("I136").Value = 3867.3333
$ws.Range("K136").Value = 11601.9999
$ws.Range("M136").Value = -9051.999899999999

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 9302.666999999999
$ws.Range("H84").Value = 9302.666999999999
$ws.Range("H113").Value = 1078.909
$ws.Range("J113").Value = 1278.8
$ws.Range("L113").Value = 3836.4
$ws.Range("N113").Value = -8176.4
